$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New timesheet entries: 11.06.22 (row 13) and 12.06.22 (row 14) ---
# Column A in this sheet stores the date as plain text (shared string,
# formatted dd/mm/yy), so force text entry (like the existing rows) by
# switching the number format to "@" before writing, then restoring the
# original dd/mm/yy display format used by the other date cells in column A.

$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = "11.06.22"
$ws.Range("A13").NumberFormat = "dd/mm/yy"
$ws.Range("B13").Value = 0.333333333333333
$ws.Range("C13").Value = 0.416666666666667
$ws.Range("D13").Formula = "=C13-B13"

$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = "12.06.22"
$ws.Range("A14").NumberFormat = "dd/mm/yy"
$ws.Range("B14").Value = 0.708333333333333
$ws.Range("C14").Value = 0.833333333333333
$ws.Range("D14").Formula = "=C14-B14"

# --- Sheet default column width drifted slightly (11.60546875 -> 11.625) ---
$ws.StandardWidth = 11.625

# --- Selection moved from E10 to E20 ---
$ws.Range("E20").Select() | Out-Null
